# Append new scraped listings to the top of the "ランサーズ" feed and
# refresh the "取得日時" timestamp for every row to 2026-01-02 12:37:43.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2026-01-02 12:37:43"

# --- 1. Make room for the 3 brand-new rows -----------------------------
# Final layout (rows 2-8):
#   2 NEW  製造業向け図面自動生成システム...
#   3 NEW  施設管理・現場業務向け チェックリスト...
#   4 (was old row 2) 【介護業務効率化】研修事業の自動化を実現するプロ募集
#   5 NEW  ホットペッパービューティーブログ一括投稿システム開発
#   6 (was old row 3) 複数WEBサイトへの日記一括投稿ツールの修正...
#   7 (was old row 4) 【報告書自動化】GASで効率的な作成フローを実現!
#   8 (was old row 5) 進行管理およびチームディレクションを担当
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(5).Insert()

# --- 2. Write the two completely new rows that land above the old data -
$ws.Range("A2").Value = $newTimestamp
$ws.Range("B2").Value = "製造業向け図面自動生成システムの開発・ツール化を支援してくださるエンジニア募集(AI/バックエンド)"
$ws.Range("C2").Value = "システム開発"
$ws.Range("D2").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E2").Value = "期限情報なし"
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5460562"
$ws.Range("G2").Value = 435
$ws.Range("H2").Value = "🔥AI,Ai ◆ツール,開発"

$ws.Range("A3").Value = $newTimestamp
$ws.Range("B3").Value = "施設管理・現場業務向け チェックリスト業務の自動化・報告書作成システム開発エンジニア募集"
$ws.Range("C3").Value = "システム開発"
$ws.Range("D3").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E3").Value = "期限情報なし"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5460563"
$ws.Range("G3").Value = 220
$ws.Range("H3").Value = "◆開発,システム開発 ◇管理"

# --- 3. The row that used to be row 2 (now row 4): refresh timestamp ---
$ws.Range("A4").Value = $newTimestamp

# --- 4. New row inserted between the old row2 and old row3 (now row 5) -
$ws.Range("A5").Value = $newTimestamp
$ws.Range("B5").Value = "ホットペッパービューティーブログ一括投稿システム開発"
$ws.Range("C5").Value = "システム開発"
$ws.Range("D5").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E5").Value = "期限情報なし"
$ws.Range("F5").Value = "https://www.lancers.jp/work/detail/5455160"
$ws.Range("G5").Value = 113
$ws.Range("H5").Value = "◆開発,システム開発"

# --- 5. The remaining old rows (now rows 6,7,8): refresh timestamp only
$ws.Range("A6").Value = $newTimestamp
$ws.Range("A7").Value = $newTimestamp
$ws.Range("A8").Value = $newTimestamp

# --- 6. Column widths ----------------------------------------------------
# ColumnWidth round-trips through OOXML with a constant +5/6 character
# offset, so subtract it here to land on the exact target widths.
$pad = 5 / 6
$ws.Columns.Item(2).ColumnWidth = 52 - $pad
$ws.Columns.Item(4).ColumnWidth = 28 - $pad
$ws.Columns.Item(8).ColumnWidth = 16 - $pad

# --- 7. Rebuild the hyperlinks on column F for all 7 data rows ---------
# Inserting rows does not renumber the sheet's existing hyperlink refs,
# so drop them all and re-add one per (now correctly positioned) row.
$ws.Cells.Item(2, 6).Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5460562", "", "", "https://www.lancers.jp/work/detail/5460562")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5460563", "", "", "https://www.lancers.jp/work/detail/5460563")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5464016", "", "", "https://www.lancers.jp/work/detail/5464016")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5455160", "", "", "https://www.lancers.jp/work/detail/5455160")
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5463948", "", "", "https://www.lancers.jp/work/detail/5463948")
$ws.Hyperlinks.Add($ws.Range("F7"), "https://www.lancers.jp/work/detail/5464025", "", "", "https://www.lancers.jp/work/detail/5464025")
$ws.Hyperlinks.Add($ws.Range("F8"), "https://www.lancers.jp/work/detail/5418064", "", "", "https://www.lancers.jp/work/detail/5418064")
